$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 847
$ws1.Range("F9").Value = 502
$ws1.Range("F16").Value = 405
$ws1.Range("F21").Value = 7515
$ws1.Range("F25").Value = 23
$ws1.Range("F26").Value = 1167
$ws1.Range("F32").Value = 203
$ws1.Range("F36").Value = 145
$ws1.Range("F40").Value = 1700

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 847
$ws4.Range("F11").Value = 502
$ws4.Range("F20").Value = 405
$ws4.Range("F25").Value = 7515
$ws4.Range("F29").Value = 23
$ws4.Range("F30").Value = 1167
$ws4.Range("F37").Value = 203
$ws4.Range("F41").Value = 145
$ws4.Range("F45").Value = 1700
